# Applies the "implemented GET batches and create user" commit to the workbook.
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Batch sheet: refresh sample batch-name data + GET-batch status codes
# ---------------------------------------------------------------------------
$batch = $wb.Worksheets.Item("Batch")

$batch.Range("C2").Value = "MLBatch3456"
$batch.Range("C3").Value = "MLBatch5667"
$batch.Range("C4").Value = "MLBatch7869"
$batch.Range("C5").Value = "MLBatch0750"

# GET batch scenarios: expected status codes swapped around
$batch.Range("G9").Value = 400
$batch.Range("D10").Value = 2
$batch.Range("G11").Value = 404

$batch.Range("C14").Value = "$%#"

# Row 17 no longer carries batch-name/class/program sample data
$batch.Range("C17:F17").ClearContents()

$batch.Range("C18").Value = "%^&"

$batch.Range("G11").Select() | Out-Null

# ---------------------------------------------------------------------------
# 2. New "User" sheet: create-user test data
# ---------------------------------------------------------------------------
$user = $wb.Worksheets.Add($null, $batch)
$user.Name = "User"

$user.Range("A1").Value = "Scenario"
$user.Range("B1").Value = "userComments"
$user.Range("C1").Value = "userEduPg"
$user.Range("D1").Value = "userEduUg"
$user.Range("E1").Value = "userFirstName"
$user.Range("F1").Value = "userLastName"
$user.Range("G1").Value = "userLinkedinUrl"
$user.Range("H1").Value = "userLocation"
$user.Range("I1").Value = "userMiddleName"
$user.Range("J1").Value = "userPhoneNumber"
$user.Range("K1").Value = "roleId"
$user.Range("L1").Value = "userRoleStatus"
$user.Range("M1").Value = "userTimeZone"
$user.Range("N1").Value = "userVisaStatus"
$user.Range("O1").Value = "loginStatus"
$user.Range("P1").Value = "userLoginEmail"
$user.Range("Q1").Value = "ExpectedStatusCode"

$user.Range("A2").Value = "CreateUserWithValidData"
$user.Range("B2").Value = "good"
$user.Range("C2").Value = "bs"
$user.Range("D2").Value = "ms"
$user.Range("E2").Value = "ram"
$user.Range("F2").Value = "wre"
$user.Range("H2").Value = "PA"
$user.Range("I2").Value = "go"
$user.Range("K2").Value = "R01"
$user.Range("L2").Value = "Active"
$user.Range("M2").Value = "EST"
$user.Range("N2").Value = "H4"
$user.Range("O2").Value = "Active"
$user.Range("P2").Value = "Ninja@gmail.com"
$user.Range("Q2").Value = 201

# Phone number keeps its leading "+" as literal text (quote-prefixed),
# matching how Excel stores a typed '+91 1236004670 entry.
$user.Range("J2").Formula = "'+91 1236004670"

# Hyperlinked LinkedIn URL -> applies the built-in Hyperlink style.
$user.Hyperlinks.Add($user.Range("G2"), "http://www.linkedin.com/", $null, $null, "www.linkedin.com")

$user.Range("J2").Select() | Out-Null

$batch.Select() | Out-Null
